$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data (A1:A12) entirely, then lay down the new A1:D4 table.
$ws.Cells.Clear()

$values = @(
    @(1, "AAA", 1, 2),
    @(2, "AAA", 3, 4),
    @(4, "AAA", 5, 2),
    @(3, "BBB", 4, 1)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Range("E11").Select() | Out-Null
